$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'64.406.01"
$ws.Range('E2').Value = '  -3.25%  '
$ws.Range('D3').Value = "'3.156.05"
$ws.Range('E3').Value = '  -2.54%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'608.32"
$ws.Range('E5').Value = '  +0.69%  '
$ws.Range('D6').Value = "'146.78"
$ws.Range('E6').Value = '  -6.40%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = "'3.153.54"
$ws.Range('E8').Value = '  -2.56%  '
$ws.Range('D9').Value = "'0.528"
$ws.Range('E9').Value = '  -3.56%  '
$ws.Range('E10').Value = '  -7.92%  '
$ws.Range('D11').Value = "'5.55"
$ws.Range('E11').Value = '  -3.70%  '
$ws.Range('D12').Value = "'0.475"
$ws.Range('E12').Value = '  -5.36%  '
$ws.Range('D13').Value = "'0.0000256"
$ws.Range('E13').Value = '  -6.13%  '
$ws.Range('D14').Value = "'36.05"
$ws.Range('E14').Value = '  -7.33%  '
$ws.Range('D15').Value = "'3.670.81"
$ws.Range('E15').Value = '  -2.42%  '
$ws.Range('D16').Value = "'64.359.73"
$ws.Range('E16').Value = '  -3.36%  '
$ws.Range('E17').Value = '  +1.02%  '
$ws.Range('D18').Value = "'3.155.50"
$ws.Range('E18').Value = '  -2.38%  '
$ws.Range('D19').Value = "'6.95"
$ws.Range('E19').Value = '  -4.68%  '
$ws.Range('D20').Value = "'478.74"
$ws.Range('E20').Value = '  -5.71%  '
$ws.Range('D21').Value = "'14.60"
$ws.Range('E21').Value = '  -4.54%  '
$ws.Range('D22').Value = "'0.710"
$ws.Range('E22').Value = '  -4.38%  '
$ws.Range('D23').Value = "'7.75"
$ws.Range('E23').Value = '  -3.21%  '
$ws.Range('D24').Value = "'13.77"
$ws.Range('E24').Value = '  -5.58%  '
$ws.Range('D25').Value = "'83.45"
$ws.Range('E25').Value = '  -3.01%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = "'2.89"
$ws.Range('E27').Value = '  -3.69%  '
$ws.Range('D28').Value = "'8.42"
$ws.Range('E28').Value = '  -6.97%  '
$ws.Range('D29').Value = "'2.19"
$ws.Range('E29').Value = '  -6.58%  '
$ws.Range('D30').Value = "'0.118"
$ws.Range('E30').Value = '  -30.41%  '
$ws.Range('D31').Value = "'6.83"
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('D32').Value = "'2.76"
$ws.Range('E32').Value = '  -5.36%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').Value = "'26.21"
$ws.Range('E34').Value = '  -6.93%  '
$ws.Range('E35').Value = '  -5.11%  '
$ws.Range('D36').Value = "'6.02"
$ws.Range('E36').Value = '  -5.24%  '
$ws.Range('E37').Value = '  -2.13%  '
$ws.Range('D38').Value = "'0.0₃0719"
$ws.Range('E38').Value = '  -10.87%  '
$ws.Range('D39').Value = "'452.64"
$ws.Range('E39').Value = '  -8.46%  '
$ws.Range('D40').Value = "'2.92"
$ws.Range('E40').Value = '  -9.19%  '
$ws.Range('D41').Value = "'0.0397"
$ws.Range('E41').Value = '  -6.00%  '
$ws.Range('D42').Value = "'8.45"
$ws.Range('E42').Value = '  -3.09%  '
$ws.Range('D43').Value = "'0.119"
$ws.Range('E43').Value = '  -7.14%  '
$ws.Range('D44').Value = "'2.844.75"
$ws.Range('E44').Value = '  -3.56%  '
$ws.Range('D45').Value = "'0.269"
$ws.Range('E45').Value = '  -7.85%  '
$ws.Range('D46').Value = "'2.27"
$ws.Range('E46').Value = '  -7.86%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = "'26.42"
$ws.Range('E47').Value = '  -6.24%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value = "'0.998"
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('D49').Value = "'2.31"
$ws.Range('E49').Value = '  -4.14%  '
$ws.Range('D50').Value = "'0.114"
$ws.Range('E50').Value = '  -3.94%  '
$ws.Range('D51').Value = "'118.57"
$ws.Range('E51').Value = '  -1.93%  '
